$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range entirely first
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Profile"
$ws.Range("B1").Value = "Account"
$ws.Range("C1").Value = "Date"

# Row 2
$ws.Range("A2").Value = "Outlook"
$ws.Range("B2").Value = "madajabbar22@gmail.com"

# Row 3 (no A3)
$ws.Range("B3").Value = "madajabbar@student.untan.ac.id"

# Row 4
$ws.Range("A4").Value = "Outlook Rpa"
$ws.Range("B4").Value = "madajabbar@student.untan.ac.id"

# Apply fill-related style to B2:B4 (matches xf applyFill=1, same visual - no actual fill color change)
$ws.Range("B2:B4").Interior.ColorIndex = -4142

# Column width (~20.43 chars, matches source workbook's A-column sizing)
$ws.Columns.Item(1).ColumnWidth = 19.6

# Selection
$ws.Range("B4").Select()
